# Normalize the "Recorded By" (column G) audit-trail lists: the authoritative
# recorder (e.g. "System") now appears first, with the former ordering of the
# remaining entries preserved in reverse — i.e. each comma-separated cell is
# simply reversed end-to-end.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)   # Column G = "Recorded By"
    $raw = $cell.Value2

    if ($raw -ne $null -and $raw -is [string] -and $raw.Contains(",")) {
        $parts = $raw -split ","
        for ($i = 0; $i -lt $parts.Length; $i++) {
            $parts[$i] = $parts[$i].Trim()
        }

        # Build the reversed list manually (index-walk) — the runtime's
        # [array]::Reverse() does not mutate PowerShell array variables here.
        $reversedParts = @()
        for ($i = $parts.Length - 1; $i -ge 0; $i--) {
            $reversedParts += $parts[$i]
        }

        $newValue = [string]::Join(", ", $reversedParts)
        $cell.Value = $newValue
    }
}
